# Intro to VR 3 - Import Models: "Updated Power Points 1 and 4"
#
# 1) Delete slide 5 ("Add some extra cubes")
# 2) Refresh the cached "datetimeFigureOut" footer field (2/11/2019 -> 2/18/2019)
#    on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- 1) Remove slide 5 -----------------------------------------------------
$p.Slides.Item(5).Delete()

# --- 2) Update the cached date text ----------------------------------------
$newDate = "2/18/2019"
$ppPlaceholderDate = 16

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $isDate = $false
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDate = $true
                }
            } catch {
                $isDate = $false
            }
            if ($isDate -eq $true) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}
